$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3-36 with shifted weekly data (D, J, K, L, M, P columns)
$ws.Range("D3").Value = 44670
$ws.Range("J3").Value = 480
$ws.Range("K3").Value = 14500
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14750
$ws.Range("P3").Value = 1135
$ws.Range("D4").Value = 44418
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14500
$ws.Range("P4").Value = 1115
$ws.Range("D5").Value = 44505
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 17000
$ws.Range("M5").Value = 16500
$ws.Range("P5").Value = 1269
$ws.Range("D6").Value = 44657
$ws.Range("J6").Value = 460
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 1192
$ws.Range("D7").Value = 44631
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 16000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 16500
$ws.Range("P7").Value = 1269
$ws.Range("D8").Value = 44344
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 18500
$ws.Range("L8").Value = 19000
$ws.Range("M8").Value = 18750
$ws.Range("P8").Value = 1442
$ws.Range("D9").Value = 44249
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 42000
$ws.Range("L9").Value = 43000
$ws.Range("M9").Value = 42500
$ws.Range("P9").Value = 3269
$ws.Range("D10").Value = 44419
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("P10").Value = 1115
$ws.Range("D11").Value = 44379
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("P11").Value = 1346
$ws.Range("D12").Value = 44414
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("P12").Value = 1115
$ws.Range("D13").Value = 44400
$ws.Range("J13").Value = 600
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("P13").Value = 1192
$ws.Range("D14").Value = 44484
$ws.Range("J14").Value = 360
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("P14").Value = 1115
$ws.Range("D15").Value = 44445
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 13000
$ws.Range("L15").Value = 14000
$ws.Range("M15").Value = 13500
$ws.Range("P15").Value = 1038
$ws.Range("D16").Value = 44658
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 15500
$ws.Range("P16").Value = 1192
$ws.Range("D17").Value = 44323
$ws.Range("J17").Value = 460
$ws.Range("K17").Value = 25000
$ws.Range("L17").Value = 26000
$ws.Range("M17").Value = 25500
$ws.Range("P17").Value = 1962
$ws.Range("D18").Value = 44333
$ws.Range("J18").Value = 440
$ws.Range("K18").Value = 24000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = 24500
$ws.Range("P18").Value = 1885
$ws.Range("D19").Value = 44533
$ws.Range("J19").Value = 520
$ws.Range("K19").Value = 17000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 17500
$ws.Range("P19").Value = 1346
$ws.Range("D20").Value = 44442
$ws.Range("J20").Value = 460
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14500
$ws.Range("P20").Value = 1115
$ws.Range("D21").Value = 44644
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 16000
$ws.Range("M21").Value = 15500
$ws.Range("P21").Value = 1192
$ws.Range("D22").Value = 44312
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 26000
$ws.Range("L22").Value = 27000
$ws.Range("M22").Value = 26500
$ws.Range("P22").Value = 2038
$ws.Range("D23").Value = 44412
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14500
$ws.Range("P23").Value = 1115
$ws.Range("D24").Value = 44365
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 19500
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 19750
$ws.Range("P24").Value = 1519
$ws.Range("D25").Value = 44575
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("P25").Value = 1115
$ws.Range("D26").Value = 44582
$ws.Range("J26").Value = 520
$ws.Range("K26").Value = 15000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 15500
$ws.Range("P26").Value = 1192
$ws.Range("D27").Value = 44596
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 16000
$ws.Range("L27").Value = 17000
$ws.Range("M27").Value = 16500
$ws.Range("P27").Value = 1269
$ws.Range("D28").Value = 44599
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 16000
$ws.Range("M28").Value = 15500
$ws.Range("P28").Value = 1192
$ws.Range("D29").Value = 44498
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 14000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 14500
$ws.Range("P29").Value = 1115
$ws.Range("D30").Value = 44326
$ws.Range("J30").Value = 460
$ws.Range("K30").Value = 25000
$ws.Range("L30").Value = 26000
$ws.Range("M30").Value = 25500
$ws.Range("P30").Value = 1962
$ws.Range("D31").Value = 44428
$ws.Range("J31").Value = 480
$ws.Range("K31").Value = 14000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 14500
$ws.Range("P31").Value = 1115
$ws.Range("D32").Value = 44309
$ws.Range("J32").Value = 400
$ws.Range("K32").Value = 26000
$ws.Range("L32").Value = 27000
$ws.Range("M32").Value = 26500
$ws.Range("P32").Value = 2038
$ws.Range("D33").Value = 44426
$ws.Range("J33").Value = 460
$ws.Range("K33").Value = 14000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 14500
$ws.Range("P33").Value = 1115
$ws.Range("D34").Value = 44335
$ws.Range("J34").Value = 480
$ws.Range("K34").Value = 24500
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = 24750
$ws.Range("P34").Value = 1904
$ws.Range("D35").Value = 44410
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 14000
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = 14500
$ws.Range("P35").Value = 1115
$ws.Range("D36").Value = 44435
$ws.Range("J36").Value = 480
$ws.Range("K36").Value = 13000
$ws.Range("L36").Value = 14000
$ws.Range("M36").Value = 13500
$ws.Range("P36").Value = 1038

# Add new row 37 (previous row 36 data, now appended as newest historical entry)
$ws.Range("A37").Value = 8
$ws.Range("B37").Value = "Terminal La Palmera de La Serena"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 44383
$ws.Range("D37").NumberFormat = $ws.Range("D36").NumberFormat()
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 100114007
$ws.Range("G37").Value = "Jengibre"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 200
$ws.Range("K37").Value = 17000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 17500
$ws.Range("N37").Value = "`$/caja 13 kilos"
$ws.Range("O37").Value = "Perú"
$ws.Range("P37").Value = 1346
$ws.Range("Q37").Value = 13
$ws.Range("R37").Value = "Hortaliza"

Write-Host "Edit complete"
